# Apply "added hindi report support" edit to the Summary workbook.
# - Updates Lower/Higher score bucket boundaries in columns B/C (rows 6-41)
# - Scrolls the sheet view back to the top (removes the stuck topLeftCell="A19")
# - Turns on AutoFilter for the data range A1:D41 (and the resulting hidden
#   _FilterDatabase defined name that Excel writes alongside it)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated Lower/Higher score values ---------------------------------
$ws.Range("C6").Value = 4

$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 8

$ws.Range("B8").Value = 9
$ws.Range("C8").Value = 12

$ws.Range("B9").Value = 13
$ws.Range("C9").Value = 16

$ws.Range("C10").Value = 2

$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 4

$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 6

$ws.Range("B13").Value = 7
$ws.Range("C13").Value = 8

$ws.Range("C14").Value = 5

$ws.Range("B15").Value = 6
$ws.Range("C15").Value = 10

$ws.Range("B16").Value = 11
$ws.Range("C16").Value = 15

$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 20

$ws.Range("C22").Value = 10

$ws.Range("B23").Value = 11
$ws.Range("C23").Value = 20

$ws.Range("B24").Value = 21
$ws.Range("C24").Value = 30

$ws.Range("B25").Value = 31
$ws.Range("C25").Value = 40

$ws.Range("C34").Value = 9

$ws.Range("B35").Value = 10
$ws.Range("C35").Value = 18

$ws.Range("B36").Value = 19
$ws.Range("C36").Value = 27

$ws.Range("B37").Value = 28
$ws.Range("C37").Value = 36

$ws.Range("C38").Value = 42

$ws.Range("B39").Value = 43
$ws.Range("C39").Value = 84

$ws.Range("B40").Value = 85
$ws.Range("C40").Value = 126

$ws.Range("B41").Value = 127
$ws.Range("C41").Value = 168

# --- Reset scroll position back to the top ------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- Turn on AutoFilter for the used range -------------------------------
$rng = $ws.Range("A1:D41")
$rng.AutoFilter()

# Excel records the filtered range as a hidden, sheet-scoped defined name
# (xlnm._FilterDatabase) when AutoFilter is applied - recreate that too.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$41")
$filterName.Visible = $false
